$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a pure number need to be forced to Text
# format before assignment (so Excel keeps them as strings, matching the
# original inlineStr cell type), then the number format is reset back to
# "General" via the Normal style so no stray style index lingers on the cell.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "515.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0928"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "492.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.296"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.116"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.08"
$ws.Range("D51").Style = "Normal"

# Cells whose new text is safe to assign directly (Excel will not
# reinterpret them as numbers: multi-dot prices, percentage strings with
# surrounding spaces/sign, or plain non-numeric text).
$ws.Range("D2").Value = "66.257.62"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "3.229.23"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +4.88%  "
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.225.28"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").Value = "3.740.54"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "66.317.12"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("E17").Value = "  +3.75%  "
$ws.Range("D18").Value = "3.221.02"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("E21").Value = "  +6.06%  "
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +4.46%  "
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("E29").Value = "  +2.93%  "
$ws.Range("E30").Value = "  +4.97%  "
$ws.Range("E31").Value = "  +9.26%  "
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("E38").Value = "  +2.64%  "
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("E41").Value = "  +2.89%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.044.83"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").Value = "0.0₃0653"
$ws.Range("E46").Value = "  +7.52%  "
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("E51").Value = "  -0.89%  "
